$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: "email" header plus three recruiter/candidate email addresses,
# matching the pattern already used in column C (hyperlinked mailto: links).
$ws.Range("E1").Value = "email"

$ws.Range("E2").Value = "rayanakhtar120330@gmail.com"
$ws.Range("E3").Value = "r.ak200330@gmail.com"
$ws.Range("E4").Value = "rayan.ak12321@gmail.com"
$ws.Range("E5").Value = "rayanakhtar120330@gmail.com"

# Add hyperlinks (mailto:) for the new email cells.
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:rayanakhtar120330@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:r.ak200330@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:rayan.ak12321@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:rayanakhtar120330@gmail.com")

# Re-apply the same "Hyperlink" cell style used by column C so the emails
# look/behave consistently (reuses the existing style rather than creating a new one).
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"

# Resize the new column to fit its (longer) email content.
$ws.Columns.Item(5).AutoFit()

# Reflect the last-selected cell in the saved view state.
$ws.Range("E9").Select()
